# Disbursement.xlsx - "Updated basic skeleton to all the tests"
#
# The ValidationTestData sheet's stipulation-related fields are renamed to a
# "*0" convention (Stipulation0/RequiredFor0/Description0/Comments0) and a
# second, duplicated block of stipulation columns (Stipulation1/RequiredFor1/
# Description1/Comments1) is appended so the sheet can carry two stipulation
# rows per test case. A handful of other header labels are also normalised
# (casing / wording) and the workbook is switched back to automatic
# calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidationTestData")

# --- Workbook-level: switch off manual calculation -------------------------
$excel.Calculation = -4105   # xlCalculationAutomatic

# --- Row 1 header label renames (same cell, new text) ----------------------
$ws.Range("B1").Value2 = "ApplicationType"
$ws.Range("D1").Value2 = "RejectApplyingFor"
$ws.Range("E1").Value2 = "ReviewApplyingFor"
$ws.Range("H1").Value2 = "SolveFor"
$ws.Range("I1").Value2 = "PaymentFrequency"
$ws.Range("J1").Value2 = "RequestedAmount"
$ws.Range("K1").Value2 = "Term"
$ws.Range("L1").Value2 = "Stipulation0"
$ws.Range("M1").Value2 = "RequiredFor0"
$ws.Range("N1").Value2 = "Description0"
$ws.Range("O1").Value2 = "Comments0"

# --- Row 3 data renames ------------------------------------------------------
$ws.Range("L3").Value2 = "Pest Inspection"
$ws.Range("M3").Value2 = "Disbursement"

# --- Append a second stipulation block in columns P:S -----------------------
# Inherit number formats / fills / borders from the existing L:O block before
# writing the new values.
$ws.Range("L1:O3").Copy()
$ws.Range("P1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P1").Value2 = "Stipulation1"
$ws.Range("Q1").Value2 = "RequiredFor1"
$ws.Range("R1").Value2 = "Description1"
$ws.Range("S1").Value2 = "Comments1"

$ws.Range("P2").Value2 = "AddStipulations"

$ws.Range("P3").Value2 = "2 years tax returns"
$ws.Range("Q3").Value2 = "Approval"
$ws.Range("R3").Value2 = "Testing"
$ws.Range("S3").Value2 = "Akcelerant"

# New columns pick up the same red "contains '='" warning conditional format
# that already exists on the rest of the sheet (NOT(ISERROR(SEARCH("=",P1)))).
$warn = $ws.Range("P1:S3").FormatConditions.Add(2, 0, '=NOT(ISERROR(SEARCH("=",P1)))')
$warn.Font.Color = 393372
$warn.Interior.Color = 13551615

# ... and the duplicate-value header highlight that every other header column
# carries (five stacked rules, same as the rest of row 1).
for ($i = 0; $i -lt 5; $i++) {
    $dup = $ws.Range("P1:S1").FormatConditions.AddUniqueValues()
    $dup.DupeUnique = 1
    $dup.Font.Color = 393372
    $dup.Interior.Color = 13551615
}

# --- Selection / view state --------------------------------------------------
$ws.Activate()
$ws.Range("E1").Select()
